$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.3
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4.33
$ws.Range("M2").Value = 1.13
$ws.Range("N2").Value = 6
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.78
$ws.Range("U2").Value = 4.6
$ws.Range("V2").Value = 1.19
$ws.Range("AD2").Value = 9.5
$ws.Range("AO2").Value = 15
$ws.Range("AR2").Value = 34

# Row 3
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9

# Row 4
$ws.Range("G4").Value = 2.63
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 3.6
$ws.Range("AO4").Value = 13

# Row 5
$ws.Range("G5").Value = 3.05
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.85
$ws.Range("L5").Value = 3.7
$ws.Range("AA5").Value = 2.32
$ws.Range("AB5").Value = 1.53
$ws.Range("AC5").Value = 5.9
$ws.Range("AD5").Value = 13.5
$ws.Range("AH5").Value = 70
$ws.Range("AL5").Value = 175
$ws.Range("AP5").Value = 11.75
$ws.Range("AQ5").Value = 40

# Row 6
$ws.Range("H6").Value = 2.65
$ws.Range("I6").Value = 2.75
$ws.Range("J6").Value = 3.7
$ws.Range("L6").Value = 3.5
$ws.Range("M6").Value = 1.15
$ws.Range("N6").Value = 4.7
$ws.Range("O6").Value = 1.62
$ws.Range("P6").Value = 2.15
$ws.Range("S6").Value = 2.82
$ws.Range("T6").Value = 1.38
$ws.Range("W6").Value = 5.2
$ws.Range("X6").Value = 1.13
$ws.Range("AA6").Value = 2.2
$ws.Range("AB6").Value = 1.6
$ws.Range("AC6").Value = 6.3
$ws.Range("AD6").Value = 13
$ws.Range("AE6").Value = 11.5
$ws.Range("AG6").Value = 35
$ws.Range("AH6").Value = 55
$ws.Range("AI6").Value = 4.7
$ws.Range("AJ6").Value = 5.4
$ws.Range("AK6").Value = 19
$ws.Range("AL6").Value = 150
$ws.Range("AN6").Value = 5.9
$ws.Range("AO6").Value = 11.75
$ws.Range("AS6").Value = 55

# Row 8
$ws.Range("S8").Value = 2.2
$ws.Range("T8").Value = 1.65

# Row 9
$ws.Range("K9").Value = 1.95
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("O9").Value = 1.44
$ws.Range("P9").Value = 2.63
$ws.Range("Q9").Value = 1.85
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2.4
$ws.Range("T9").Value = 1.53
$ws.Range("AC9").Value = 6.5
$ws.Range("AN9").Value = 7.5

# Row 11
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("Q11").Value = 1.83
$ws.Range("R11").Value = 2.03
$ws.Range("S11").Value = 2.4
$ws.Range("T11").Value = 1.53
$ws.Range("W11").Value = 4.5
$ws.Range("X11").Value = 1.18

# Row 12
$ws.Range("H12").Value = 3.2
$ws.Range("I12").Value = 3.3
$ws.Range("J12").Value = 3
$ws.Range("K12").Value = 1.95
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("Q12").Value = 1.83
$ws.Range("R12").Value = 2.03
$ws.Range("W12").Value = 4.5
$ws.Range("X12").Value = 1.18
$ws.Range("Y12").Value = 1.53
$ws.Range("Z12").Value = 2.38
$ws.Range("AA12").Value = 2.05
$ws.Range("AB12").Value = 1.7
$ws.Range("AE12").Value = 9.5
$ws.Range("AH12").Value = 34
$ws.Range("AI12").Value = 7

# Row 14
$ws.Range("G14").Value = 2.2
$ws.Range("H14").Value = 3
$ws.Range("I14").Value = 3.4
$ws.Range("J14").Value = 2.88
$ws.Range("K14").Value = 1.95
$ws.Range("L14").Value = 4
$ws.Range("M14").Value = 1.08
$ws.Range("N14").Value = 7
$ws.Range("O14").Value = 1.4
$ws.Range("P14").Value = 2.75
$ws.Range("S14").Value = 2.3
$ws.Range("T14").Value = 1.6
$ws.Range("W14").Value = 4.33
$ws.Range("X14").Value = 1.2
$ws.Range("Y14").Value = 1.53
$ws.Range("Z14").Value = 2.38
$ws.Range("AA14").Value = 2
$ws.Range("AB14").Value = 1.73
$ws.Range("AC14").Value = 6.5
$ws.Range("AD14").Value = 9.5
$ws.Range("AE14").Value = 10
$ws.Range("AF14").Value = 21
$ws.Range("AG14").Value = 21
$ws.Range("AH14").Value = 34
$ws.Range("AI14").Value = 7
$ws.Range("AJ14").Value = 6
$ws.Range("AK14").Value = 17
$ws.Range("AM14").Value = 900
$ws.Range("AN14").Value = 9
$ws.Range("AO14").Value = 15
$ws.Range("AP14").Value = 13
$ws.Range("AR14").Value = 29
$ws.Range("AS14").Value = 41
